$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8 and 9 (the "Neutrophils" sending-cluster rows)
$ws.Rows.Item(8).EntireRow.Delete()
$ws.Rows.Item(8).EntireRow.Delete()

# Update rows 2-7 with the new values
# Row 2: ECs, Wnt9a, Fzd10, FAPs
$ws.Range("G2").Value = 0.5939095
$ws.Range("H2").Value = 1.187819
$ws.Range("I2").Value = 0.4800062555800474
$ws.Range("J2").Value = 0.4083203479364477
$ws.Range("O2").Value = 0.4086672402490986
$ws.Range("P2").Value = 0.5089958879585649
$ws.Range("Q2").Value = 0.0236963951405
$ws.Range("R2").Value = 0.142178370843
$ws.Range("S2").Value = 0.1961628317702015
$ws.Range("T2").Value = 0.2078333780694623

# Row 3: ECs, Wnt9a, Fzd10, MuSCs
$ws.Range("G3").Value = 0.5939095
$ws.Range("H3").Value = 1.187819
$ws.Range("I3").Value = 0.4800062555800474
$ws.Range("J3").Value = 0.4083203479364477
$ws.Range("M3").Value = 0.057733
$ws.Range("N3").Value = 0.115466
$ws.Range("O3").Value = 0.5913327597509014
$ws.Range("P3").Value = 0.4910041120414351
$ws.Range("Q3").Value = 0.0342881771635
$ws.Range("R3").Value = 0.137152708654
$ws.Range("S3").Value = 0.283843423809846
$ws.Range("T3").Value = 0.2004869698669853

# Row 4: FAPs, Wnt9a, Fzd10, FAPs
$ws.Range("I4").Value = 0.3511258158251627
$ws.Range("J4").Value = 0.4480310838260221
$ws.Range("O4").Value = 0.4086672402490986
$ws.Range("P4").Value = 0.5089958879585649
$ws.Range("S4").Value = 0.1434936181334825
$ws.Range("T4").Value = 0.2280459793450643

# Row 5: FAPs, Wnt9a, Fzd10, MuSCs
$ws.Range("I5").Value = 0.3511258158251627
$ws.Range("J5").Value = 0.4480310838260221
$ws.Range("M5").Value = 0.057733
$ws.Range("N5").Value = 0.115466
$ws.Range("O5").Value = 0.5913327597509014
$ws.Range("P5").Value = 0.4910041120414351
$ws.Range("Q5").Value = 0.02508189016233333
$ws.Range("R5").Value = 0.150491340974
$ws.Range("S5").Value = 0.2076321976916802
$ws.Range("T5").Value = 0.2199851044809578

# Row 6: MuSCs, Wnt9a, Fzd10, FAPs
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.2089395
$ws.Range("H6").Value = 0.417879
$ws.Range("I6").Value = 0.1688679285947898
$ws.Range("J6").Value = 0.1436485682375301
$ws.Range("O6").Value = 0.4086672402490986
$ws.Range("P6").Value = 0.5089958879585649
$ws.Range("Q6").Value = 0.0083364771105
$ws.Range("R6").Value = 0.050018862663
$ws.Range("S6").Value = 0.06901079034541459
$ws.Range("T6").Value = 0.07311653054403815

# Row 7: MuSCs, Wnt9a, Fzd10, MuSCs
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.2089395
$ws.Range("H7").Value = 0.417879
$ws.Range("I7").Value = 0.1688679285947898
$ws.Range("J7").Value = 0.1436485682375301
$ws.Range("M7").Value = 0.057733
$ws.Range("N7").Value = 0.115466
$ws.Range("O7").Value = 0.5913327597509014
$ws.Range("P7").Value = 0.4910041120414351
$ws.Range("Q7").Value = 0.0120627041535
$ws.Range("R7").Value = 0.048250816614
$ws.Range("S7").Value = 0.09985713824937521
$ws.Range("T7").Value = 0.07053203769349198
